$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Helper: write a value while keeping the cell's original style/format
# (the source cells hold plain text, even for numeric-looking strings like
# "244.46" - setting .Value directly on such a string gets auto-coerced to
# a Number by Excel, so we briefly force a Text number format, assign the
# value, then restore the original style so no stray formatting leaks in).
function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

# Simple price/label updates (rows unaffected by reordering)
Set-TextValue $ws.Cells.Item(2,4) "244.43"
Set-TextValue $ws.Cells.Item(3,4) "21.97"
Set-TextValue $ws.Cells.Item(4,4) "5.399"
Set-TextValue $ws.Cells.Item(6,4) "3.386"
Set-TextValue $ws.Cells.Item(7,4) "0.8117"
Set-TextValue $ws.Cells.Item(8,4) "0.9552"
Set-TextValue $ws.Cells.Item(9,4) "0.1422"
Set-TextValue $ws.Cells.Item(10,4) "0.07438"
Set-TextValue $ws.Cells.Item(11,4) "0.03366"
Set-TextValue $ws.Cells.Item(12,4) "0.03052"
Set-TextValue $ws.Cells.Item(13,4) "0.09416"
Set-TextValue $ws.Cells.Item(14,4) "4.003"
Set-TextValue $ws.Cells.Item(15,4) "0.001588"
Set-TextValue $ws.Cells.Item(16,4) "0.04831"
Set-TextValue $ws.Cells.Item(17,4) "0.0005870"
Set-TextValue $ws.Cells.Item(17,5) "16OneONEWorstin24h"
Set-TextValue $ws.Cells.Item(18,4) "0.006120"
Set-TextValue $ws.Cells.Item(20,4) "0.0009884"
Set-TextValue $ws.Cells.Item(22,4) "3.694"
Set-TextValue $ws.Cells.Item(23,4) "6.404"
Set-TextValue $ws.Cells.Item(24,4) "2.185"
Set-TextValue $ws.Cells.Item(26,4) "0.1292"
Set-TextValue $ws.Cells.Item(27,4) "0.0002460"

# Rows 41-43: coins shuffled (BKEXToken/CEJI/KickToken cycled), with new prices
Set-TextValue $ws.Cells.Item(41,2) "KickToken"
Set-TextValue $ws.Cells.Item(41,3) "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws.Cells.Item(41,4) "0.006553"
Set-TextValue $ws.Cells.Item(41,5) "40KickTokenKICK"

Set-TextValue $ws.Cells.Item(42,2) "BKEXToken"
Set-TextValue $ws.Cells.Item(42,3) "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws.Cells.Item(42,4) "0.1074"
Set-TextValue $ws.Cells.Item(42,5) "41BKEXTokenBKK"

Set-TextValue $ws.Cells.Item(43,2) "CEJI"
Set-TextValue $ws.Cells.Item(43,3) "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws.Cells.Item(43,4) "0.002901"
Set-TextValue $ws.Cells.Item(43,5) "42CEJICEJI"

Set-TextValue $ws.Cells.Item(44,4) "0.005228"
Set-TextValue $ws.Cells.Item(45,4) "0.00005413"
Set-TextValue $ws.Cells.Item(47,4) "0.8100"
Set-TextValue $ws.Cells.Item(48,4) "0.01946"
Set-TextValue $ws.Cells.Item(48,5) "47BOLOBOLO"
